$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row before the current blank separator row (row 16), shifting
# everything below down by one. The new row picks up formatting from row 15
# above it (style s="4" on column A, s="5" on column H), matching the target.
$ws.Range("A16").EntireRow.Insert()

# The worksheet's hyperlinks don't automatically track the row shift, so
# rebuild them at their new locations (B98 -> B99; H7 is unaffected).
$ws.Hyperlinks.Delete()
$null = $ws.Hyperlinks.Add($ws.Range("B99"), "http://www.officekb.com/Uwe/Forum.aspx/excel-prog/159706/Shape-TextEffect-HorizontalAlignment-throws-error")
$null = $ws.Hyperlinks.Add($ws.Range("H7"), "http://www.vbforums.com/archive/index.php/t-47843.html")
# Re-adding a hyperlink recolours the cell; restore the original Hyperlink style.
$ws.Range("B99").Style = "Hyperlink"
$ws.Range("H7").Style = "Hyperlink"

# Re-populate the bullet list for the "Version 2.5 alpha" section (rows 11-16)
# with the updated wording / ordering from the changelog entry.
$ws.Cells.Item(11, 1).Value2 = "Updated CBC.exe to version 2.8.8"
$ws.Cells.Item(12, 1).Value2 = "Support for using the Gurobi LP/IP solver if a user has this installed on their machine"
$ws.Cells.Item(13, 1).Value2 = "Support for cloud-based NEOS server for CBC solver"
$ws.Cells.Item(14, 1).Value2 = "Support for solving non-linear models using both NOMAD and the cloud-based NEOS servers (assuming non-negativity currently doesn't work correctly for non-linear NEOS, all variables are assumed positive, not just unconstrained ones)"
$ws.Cells.Item(15, 1).Value2 = "Reporting of dual variables and sensitivity analysis"
$ws.Cells.Item(16, 1).Value2 = "Many small bux fixes and feature enhancements"

# The Solver add-in's hidden "solver_opt" defined name tracks a cell reference
# that shifts down one row because of the inserted row above it.
$wb.Names.Item("Sheet1!solver_opt").RefersTo = "=Sheet1!`$U`$91"

# Switch the workbook back to automatic calculation.
$excel.Calculation = -4105

# Restore the on-screen selection to where the editor was last working.
$null = $ws.Activate()
$null = $ws.Range("I16").Select()
